# Edit script: tweaks to editors and params fixing cathedral which sounded
# terrible; cathedral & echo now use wet & dry mixes rather than wet-dry mix.
#
# This appends a large batch of new size-tracking measurements to Sheet2,
# adds a "diff" column header, and extends the running B-minus-previous-B
# delta formula down through the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---- Header: add a "diff" column title in C1 ----
$ws.Range("C1").Value = "diff"

# ---- New data rows (19 through 56) ----
# Each entry: row, A (label, optional), B (size value), D (note, optional)
$rows = @(
    @{r=19; A="all devices with M7 style param"; B=21276},
    @{r=20; A="remove rotation from width"; B=21240},
    @{r=21; A="adding simpler overload of getfrequency()"; B=21224},
    @{r=22; A="simpler overload of Get01value"; B=21220},
    @{r=23; A="w/o sat."; B=19976; D="sat is 1.2kb; unacceptable"},
    @{r=24; A="simpler overload of Get01value"; B=21220},
    @{r=25; A="no initialize cathedral values, some tweaks"; B=21180},
    @{r=26; A="no maj7 obsolete + var tri waveforms"; B=21100},
    @{r=27; A="simpler sine clip"; B=21000; D="wow that was a lot of savings :x"},
    @{r=28; A="simpler saw"; B=20940},
    @{r=29; A="w/o sat."; B=19696},
    @{r=30; A="with all devices"; B=20940},
    @{r=31; A="tiny sat optimization moving a multiply"; B=20928},
    @{r=32; A="disable 48db crossover"; B=20856},
    @{r=33; A="removing rarely used sat models"; B=20680},
    @{r=34; A="removing analog support"; B=20616},
    @{r=35; A="bypassing sat processample"; B=20344},
    @{r=36; B=20616},
    @{r=37; A="bypass distort()"; B=20492},
    @{r=38; B=20616},
    @{r=39; A="removing all div style models"; B=20572},
    @{r=40; A="different method of stereo proce"; B=20584},
    @{r=41; A="no mute/solo processing when no selectable stream"; B=20572},
    @{r=42; A="tiny fix"; B=20576},
    @{r=43; A="no ms in sat"; B=20536},
    @{r=44; A="no sat"; B=19696; D="sat is now 840 bytes"},
    @{r=45; B=20536},
    @{r=46; A="no cathedral"; B=19812; D="cathedral is 724 bytes"},
    @{r=47; B=20536},
    @{r=48; A='"optimised" echo processing'; B=20560; D="well I tried but cannot get this to be smaller."},
    @{r=49; A="inlining the simplest paramaccessors"; B=20644; D="wow that is bad."},
    @{r=50; A="echo baseline again"; B=20536},
    @{r=51; A="absolute best echo optimization; nope."; B=20548},
    @{r=52; A="baseline."; B=20536},
    @{r=53; A="unifying loaddefaults and get/setparam"; B=20468},
    @{r=54; A="some fixes after testing"; B=20524},
    @{r=55; A="comp no full features"; B=20380; D="worth it"},
    @{r=56; A="what was our baseline again"; B=21276}
)

foreach ($row in $rows) {
    $r = $row.r
    if ($row.ContainsKey("A")) {
        $ws.Range("A$r").Value = $row.A
    }
    $ws.Range("B$r").Value = $row.B
    if ($row.ContainsKey("D")) {
        $ws.Range("D$r").Value = $row.D
    }
}

# Row 19: A19, B19 and C19 all carry the yellow-highlight style (style index 2
# in styles.xml -> solid yellow fill). Row 56: only B56 carries that style.
$ws.Range("A19:C19").Interior.Color = 65535
$ws.Range("B56").Interior.Color = 65535

# ---- Extend the running delta formula (C column) down through row 73 ----
# Existing shared formula already covers C8:C18 (=B-prev B). Continue the
# same relative pattern for the newly added rows and the trailing zero rows.
for ($r = 19; $r -le 73; $r++) {
    $ws.Range("C$r").FormulaR1C1 = "=R[0]C[-1]-R[-1]C[-1]"
}

# ---- View state: mirror the author's scroll position / selection ----
$ws.Application.ActiveWindow.ScrollRow = 35
$ws.Range("D56").Select()
